$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 holds the nutrient-name header, duplicated across a "male" column
# and a "female" column for every nutrient (odd columns = male, per the
# 性別/男/女 labels in row 4). Prefix the male (odd) column headers with
# "男" so they read distinctly from their female counterparts.
$ws.Range("C2").Value  = "男維生素A"
$ws.Range("E2").Value  = "男維生素D"
$ws.Range("G2").Value  = "男維生素E"
$ws.Range("I2").Value  = "男維生素K"
$ws.Range("K2").Value  = "男維生素C"
$ws.Range("M2").Value  = "男維生素B1"
$ws.Range("O2").Value  = "男維生素B2"
$ws.Range("Q2").Value  = "男菸鹼素"
$ws.Range("S2").Value  = "男維生素B6"
$ws.Range("U2").Value  = "男維生素B12"
$ws.Range("W2").Value  = "男葉酸"
$ws.Range("Y2").Value  = "男鈣"
$ws.Range("AA2").Value = "男磷"
$ws.Range("AC2").Value = "男鎂"
$ws.Range("AE2").Value = "男鐵"
$ws.Range("AG2").Value = "男鋅"
$ws.Range("AI2").Value = "男碘"
$ws.Range("AK2").Value = "男鉀"
$ws.Range("AM2").Value = "男鈉"

# These same "male" header cells also pick up a distinct font treatment
# (a fresh bold Arial face) so they stand out visually from the rest of
# row 2.
$maleHeaderCells = @("C2","M2","O2","Q2","S2","U2","W2","Y2","AA2","AC2","AE2","AG2","AI2","AK2","AM2")
foreach ($cellRef in $maleHeaderCells) {
    $cell = $ws.Range($cellRef)
    $cell.Font.Bold = $true
    $cell.Font.Name = "Arial"
}

# Leave the view scrolled back to the left with H19 selected.
$ws.Range("H19").Select() | Out-Null
